# "Generate Report for Handback"
#
# The localization-status report is regenerated: rows get re-sorted
# alphabetically by source-file GUID (3b4240e3... and 5b03f02d... move up
# from positions 4/5 to 1/2), and the two files that were previously only
# "Ready for handoff" are now shown as handed back - same as the other two
# files already were - with their handback file / datetime columns filled
# in.
#
# Only the cells whose value actually differs between the old and new
# report are touched; everything else (column D/E/F/L/M/N/O/P, which are
# identical across all four rows both before and after) is left alone.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Hyperlinks.Delete()

$ov.Range("A2").Value = "3b4240e3-eb42-4252-9429-4e0ea964a0b9.md"
$ov.Range("B2").Value = "e2e\3b4240e3-eb42-4252-9429-4e0ea964a0b9.md"
$ov.Range("G2").Value = "2016-09-02 18:27:06"

$ov.Range("A3").Value = "5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md"
$ov.Range("B3").Value = "e2e\5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md"
$ov.Range("G3").Value = "2016-09-02 18:27:06"

$ov.Range("A4").Value = "97a8312f-140e-4cb7-811b-3ecfc064aee0.md"
$ov.Range("B4").Value = "e2e\97a8312f-140e-4cb7-811b-3ecfc064aee0.md"
$ov.Range("E4").Value = "Handed back: in sync with en-US"
$ov.Range("F4").Value = "Handed back: in sync with en-US"
$ov.Range("G4").Value = "2016-09-02 18:25:20"

$ov.Range("A5").Value = "fc294280-e350-4ba2-8853-74731166f92f.md"
$ov.Range("B5").Value = "e2e\fc294280-e350-4ba2-8853-74731166f92f.md"
$ov.Range("E5").Value = "Handed back: in sync with en-US"
$ov.Range("F5").Value = "Handed back: in sync with en-US"
$ov.Range("G5").Value = "2016-09-02 18:25:20"

$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60eaf0d34c80879aae3606f376958d761633fa3f/e2e/3b4240e3-eb42-4252-9429-4e0ea964a0b9.md", "", "", "e2e\3b4240e3-eb42-4252-9429-4e0ea964a0b9.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60eaf0d34c80879aae3606f376958d761633fa3f/e2e/5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md", "", "", "e2e\5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a2de00bc04b6e1f5c5e3f604a28b5e530eb7d74/e2e/97a8312f-140e-4cb7-811b-3ecfc064aee0.md", "", "", "e2e\97a8312f-140e-4cb7-811b-3ecfc064aee0.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a2de00bc04b6e1f5c5e3f604a28b5e530eb7d74/e2e/fc294280-e350-4ba2-8853-74731166f92f.md", "", "", "e2e\fc294280-e350-4ba2-8853-74731166f92f.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Delete()

$zh.Range("A2").Value = "3b4240e3-eb42-4252-9429-4e0ea964a0b9.md"
$zh.Range("G2").Value = "3b4240e3-eb42-4252-9429-4e0ea964a0b9.b774fcbaf8d047461330a80ac93839c050da39a0.zh-cn.xlf"
$zh.Range("H2").Value = "2016-09-02 18:26:58"
$zh.Range("I2").Value = "3b4240e3-eb42-4252-9429-4e0ea964a0b9.md"
$zh.Range("J2").Value = "3b4240e3-eb42-4252-9429-4e0ea964a0b9.b774fcbaf8d047461330a80ac93839c050da39a0.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-02 18:27:30"

$zh.Range("A3").Value = "5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md"
$zh.Range("G3").Value = "5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.43ae2c6352ced74ff3fba1f420979a435aac6045.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-02 18:26:58"
$zh.Range("I3").Value = "5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md"
$zh.Range("J3").Value = "5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.43ae2c6352ced74ff3fba1f420979a435aac6045.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-02 18:27:30"

$zh.Range("A4").Value = "97a8312f-140e-4cb7-811b-3ecfc064aee0.md"
$zh.Range("C4").Value = "Handed back: in sync with en-US"
$zh.Range("G4").Value = "97a8312f-140e-4cb7-811b-3ecfc064aee0.3d4fbe1120cb785540a9de8abb6383e9d88924df.zh-cn.xlf"
$zh.Range("H4").Value = "2016-09-02 18:25:15"
$zh.Range("I4").Value = "97a8312f-140e-4cb7-811b-3ecfc064aee0.md"
$zh.Range("J4").Value = "97a8312f-140e-4cb7-811b-3ecfc064aee0.3d4fbe1120cb785540a9de8abb6383e9d88924df.zh-cn.xlf"
$zh.Range("K4").Value = "2016-09-02 18:26:17"

$zh.Range("A5").Value = "fc294280-e350-4ba2-8853-74731166f92f.md"
$zh.Range("C5").Value = "Handed back: in sync with en-US"
$zh.Range("G5").Value = "fc294280-e350-4ba2-8853-74731166f92f.bdb0c87f0a9cff20a88920094cd4c6231b2ffecf.zh-cn.xlf"
$zh.Range("H5").Value = "2016-09-02 18:25:15"
$zh.Range("I5").Value = "fc294280-e350-4ba2-8853-74731166f92f.md"
$zh.Range("J5").Value = "fc294280-e350-4ba2-8853-74731166f92f.bdb0c87f0a9cff20a88920094cd4c6231b2ffecf.zh-cn.xlf"
$zh.Range("K5").Value = "2016-09-02 18:26:17"

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60eaf0d34c80879aae3606f376958d761633fa3f/e2e/3b4240e3-eb42-4252-9429-4e0ea964a0b9.md", "", "", "3b4240e3-eb42-4252-9429-4e0ea964a0b9.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c5cbbee5ca34c5fb8cd5c5ef9f1f09b70ec97ecc/e2e/3b4240e3-eb42-4252-9429-4e0ea964a0b9.md", "", "", "3b4240e3-eb42-4252-9429-4e0ea964a0b9.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60eaf0d34c80879aae3606f376958d761633fa3f/e2e/5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md", "", "", "5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c5cbbee5ca34c5fb8cd5c5ef9f1f09b70ec97ecc/e2e/5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md", "", "", "5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a2de00bc04b6e1f5c5e3f604a28b5e530eb7d74/e2e/97a8312f-140e-4cb7-811b-3ecfc064aee0.md", "", "", "97a8312f-140e-4cb7-811b-3ecfc064aee0.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c5cbbee5ca34c5fb8cd5c5ef9f1f09b70ec97ecc/e2e/97a8312f-140e-4cb7-811b-3ecfc064aee0.md", "", "", "97a8312f-140e-4cb7-811b-3ecfc064aee0.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a2de00bc04b6e1f5c5e3f604a28b5e530eb7d74/e2e/fc294280-e350-4ba2-8853-74731166f92f.md", "", "", "fc294280-e350-4ba2-8853-74731166f92f.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c5cbbee5ca34c5fb8cd5c5ef9f1f09b70ec97ecc/e2e/fc294280-e350-4ba2-8853-74731166f92f.md", "", "", "fc294280-e350-4ba2-8853-74731166f92f.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Delete()

$de.Range("A2").Value = "3b4240e3-eb42-4252-9429-4e0ea964a0b9.md"
$de.Range("G2").Value = "3b4240e3-eb42-4252-9429-4e0ea964a0b9.b774fcbaf8d047461330a80ac93839c050da39a0.de-de.xlf"
$de.Range("H2").Value = "2016-09-02 18:27:06"
$de.Range("I2").Value = "3b4240e3-eb42-4252-9429-4e0ea964a0b9.md"
$de.Range("J2").Value = "3b4240e3-eb42-4252-9429-4e0ea964a0b9.b774fcbaf8d047461330a80ac93839c050da39a0.de-de.xlf"
$de.Range("K2").Value = "2016-09-02 18:27:37"

$de.Range("A3").Value = "5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md"
$de.Range("G3").Value = "5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.43ae2c6352ced74ff3fba1f420979a435aac6045.de-de.xlf"
$de.Range("H3").Value = "2016-09-02 18:27:06"
$de.Range("I3").Value = "5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md"
$de.Range("J3").Value = "5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.43ae2c6352ced74ff3fba1f420979a435aac6045.de-de.xlf"
$de.Range("K3").Value = "2016-09-02 18:27:37"

$de.Range("A4").Value = "97a8312f-140e-4cb7-811b-3ecfc064aee0.md"
$de.Range("C4").Value = "Handed back: in sync with en-US"
$de.Range("G4").Value = "97a8312f-140e-4cb7-811b-3ecfc064aee0.3d4fbe1120cb785540a9de8abb6383e9d88924df.de-de.xlf"
$de.Range("H4").Value = "2016-09-02 18:25:20"
$de.Range("I4").Value = "97a8312f-140e-4cb7-811b-3ecfc064aee0.md"
$de.Range("J4").Value = "97a8312f-140e-4cb7-811b-3ecfc064aee0.3d4fbe1120cb785540a9de8abb6383e9d88924df.de-de.xlf"
$de.Range("K4").Value = "2016-09-02 18:26:24"

$de.Range("A5").Value = "fc294280-e350-4ba2-8853-74731166f92f.md"
$de.Range("C5").Value = "Handed back: in sync with en-US"
$de.Range("G5").Value = "fc294280-e350-4ba2-8853-74731166f92f.bdb0c87f0a9cff20a88920094cd4c6231b2ffecf.de-de.xlf"
$de.Range("H5").Value = "2016-09-02 18:25:20"
$de.Range("I5").Value = "fc294280-e350-4ba2-8853-74731166f92f.md"
$de.Range("J5").Value = "fc294280-e350-4ba2-8853-74731166f92f.bdb0c87f0a9cff20a88920094cd4c6231b2ffecf.de-de.xlf"
$de.Range("K5").Value = "2016-09-02 18:26:24"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60eaf0d34c80879aae3606f376958d761633fa3f/e2e/3b4240e3-eb42-4252-9429-4e0ea964a0b9.md", "", "", "3b4240e3-eb42-4252-9429-4e0ea964a0b9.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/fddc86f70a694ad7f457a9b68bf5496b2c15241f/e2e/3b4240e3-eb42-4252-9429-4e0ea964a0b9.md", "", "", "3b4240e3-eb42-4252-9429-4e0ea964a0b9.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60eaf0d34c80879aae3606f376958d761633fa3f/e2e/5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md", "", "", "5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/fddc86f70a694ad7f457a9b68bf5496b2c15241f/e2e/5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md", "", "", "5b03f02d-7b97-47a8-b51b-d7435dc5d3f4.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a2de00bc04b6e1f5c5e3f604a28b5e530eb7d74/e2e/97a8312f-140e-4cb7-811b-3ecfc064aee0.md", "", "", "97a8312f-140e-4cb7-811b-3ecfc064aee0.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/fddc86f70a694ad7f457a9b68bf5496b2c15241f/e2e/97a8312f-140e-4cb7-811b-3ecfc064aee0.md", "", "", "97a8312f-140e-4cb7-811b-3ecfc064aee0.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a2de00bc04b6e1f5c5e3f604a28b5e530eb7d74/e2e/fc294280-e350-4ba2-8853-74731166f92f.md", "", "", "fc294280-e350-4ba2-8853-74731166f92f.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/fddc86f70a694ad7f457a9b68bf5496b2c15241f/e2e/fc294280-e350-4ba2-8853-74731166f92f.md", "", "", "fc294280-e350-4ba2-8853-74731166f92f.md") | Out-Null
